# Apply the Sun Jul 16 18:00:09 UTC 2023 GitHub Actions cryptos-list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.360.51"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.935.80"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'0.7740"
$ws.Range("E5").Value = "  +8.44%  "
$ws.Range("D6").Value = "'246.66"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("D7").Value = "'0.9987"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'28.02"
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").Value = "'0.3215"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").Value = "'0.07094"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").Value = "'0.7838"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("D12").Value = "'0.08021"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "1.935.96"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'5.383"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "'95.09"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "'14.57"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "30.352.98"
$ws.Range("D18").Value = "'255.99"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'0.000008019"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").Value = "'5.826"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "2.190.12"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'0.9983"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "'6.764"
$ws.Range("E24").Value = "  -2.84%  "
$ws.Range("D25").Value = "'9.617"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").Value = "'164.09"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").Value = "'0.1351"
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("D28").Value = "'19.12"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("D30").Value = "'1.367"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "'1.521"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").Value = "'4.436"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "'4.147"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "'0.05196"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'1.288"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").Value = "'0.7540"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "'2.772"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "'0.01969"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Value = "'2.811"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "'79.14"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "'6.491"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").Value = "'0.4527"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "'1.982"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'0.8357"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "'101.25"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "'9.784"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "'7.516"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'984.69"
$ws.Range("E49").Value = "  +11.11%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'37.38"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "'0.4175"
$ws.Range("E51").Value = "  +0.40%  "
